# Insert four new "results" paragraphs immediately before the existing
# first paragraph ("Cheeseburgers and fries."), per the commit's intent
# of laying the results out (author/location/date/end) ahead of the
# existing body text.

$d = $word.ActiveDocument

# Anchor on the current first paragraph and repeatedly insert new
# paragraphs directly before it. Because each InsertBefore lands right
# at the (stable) start of that anchor paragraph, inserting in reverse
# order yields the correct final top-to-bottom sequence.
$anchor = $d.Paragraphs(1).Range

$anchor.InsertBefore("end`r")
$anchor.InsertBefore("Date May 22`r")
$anchor.InsertBefore("Location: Las Vegas`r")
$anchor.InsertBefore("Authour: Barry`r")
